$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells AD1, AE1, AF1 matching the style of existing header cells (e.g. A1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the values (PasteSpecial formats only, but ensure values remain correct)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 1
}
